# PM10 Tidsregistrering for Nikolaj - add new time-tracking rows 9-14
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 9 - Operations kontrakt OC0801 / System Analyst
$ws.Range("A9").Value = "Operations kontrakt OC0801"
$ws.Range("B9").Value = "System Analyst "
$ws.Range("C9").Value = 43887
$ws.Range("D9").Value = 0.39583333333333331
$ws.Range("E9").Value = 0.41666666666666669
$ws.Range("F9").Value = "30 min"

# Row 10 - Kravworkshop / Requirement Specifier
$ws.Range("A10").Value = "Kravworkshop"
$ws.Range("B10").Value = "Requirement Specifier"
$ws.Range("C10").Value = 43887
$ws.Range("D10").Value = 0.4375
$ws.Range("E10").Value = 0.5
$ws.Range("F10").Value = "1:30 timer"

# Row 11 - Kravworkshop / Requirement Specifier (new number-format style on F11)
$ws.Range("A11").Value = "Kravworkshop"
$ws.Range("B11").Value = "Requirement Specifier"
$ws.Range("C11").Value = 43887
$ws.Range("D11").Value = 0.52083333333333337
$ws.Range("E11").Value = 0.5625
$ws.Range("F11").Value = "2 timer og 30 min"
$ws.Range("F11").HorizontalAlignment = -4108
$ws.Range("F11").VerticalAlignment = -4108
$ws.Range("F11").NumberFormat = "h:mm"

# Row 12 - Usecase 07 indtjeningsbidrag / System Analyst
$ws.Range("A12").Value = "Usecase 07 indtjeningsbidrag"
$ws.Range("B12").Value = "System Analyst "
$ws.Range("C12").Value = 43887
$ws.Range("D12").Value = 0.5625
$ws.Range("E12").Value = 0.58333333333333337
$ws.Range("F12").Value = "1 time"

# Row 13 - Domæne model 07 (no role)
$ws.Range("A13").Value = "Domæne model 07"
$ws.Range("C13").Value = 43887
$ws.Range("D13").Value = 0.58333333333333337
$ws.Range("E13").Value = 0.60416666666666663
$ws.Range("F13").Value = "30 min"

# Row 14 - Networking (Virksomhedsdag datamatiker) (no role, F14 left blank)
$ws.Range("A14").Value = "Networking (Virksomhedsdag datamatiker)"
$ws.Range("C14").Value = 43887
$ws.Range("D14").Value = 0.60416666666666663
$ws.Range("E14").Value = 0.64583333333333337

# Update sheet view: zoom, scroll position and active selection
$ws.Activate()
$win = $excel.ActiveWindow
$win.Zoom = 119
$win.ScrollColumn = 3
$win.ScrollRow = 1
$ws.Range("F14").Select()

Write-Host "Applied Nikolaj tidsregistrering updates"
